$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted into the "Zanahoria" (Los Lagos /
# Feria Lagunitas de Puerto Montt) sheet just before the existing row 212,
# pushing every subsequent record down by one row (old row 262 -> new row 263).
$ws.Rows.Item(212).Insert()

# Populate the newly inserted row with the new observation's data.
$row = 212
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44543
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100114013
$ws.Cells.Item($row, 7).Value = "Zanahoria"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 150
$ws.Cells.Item($row, 11).Value = 10000
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 10000
$ws.Cells.Item($row, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item($row, 15).Value = "Región de Ñuble"
$ws.Cells.Item($row, 16).Value = 500
$ws.Cells.Item($row, 17).Value = 20
$ws.Cells.Item($row, 18).Value = "Hortaliza"
